$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 625.2308
$ws.Range("I103").Value = 764.3333
$ws.Range("K103").Value = 2292.9999
$ws.Range("M103").Value = -1706.9999
$ws.Range("H106").Value = 4980.75
$ws.Range("J106").Value = 4997
$ws.Range("L106").Value = 4997
$ws.Range("N106").Value = -6259
$ws.Range("H137").Value = 61909.566
$ws.Range("I137").Value = 95596.69
$ws.Range("J137").Value = 3722.7273
$ws.Range("K137").Value = 286790.07
$ws.Range("L137").Value = 11168.1819
$ws.Range("M137").Value = -284240.07
$ws.Range("N137").Value = -16268.1819
$ws.Range("H138").Value = 2932.739
$ws.Range("J138").Value = 3293.2156
$ws.Range("L138").Value = 9879.6468
$ws.Range("N138").Value = -20159.6468
$ws.Range("H141").Value = 2933.5715
$ws.Range("I141").Value = 2797
$ws.Range("K141").Value = 8391
$ws.Range("M141").Value = -3211
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3959.411
$ws.Range("I32").Value = 1833.5333
$ws.Range("K32").Value = 1833.5333
$ws.Range("M32").Value = -1546.5333
$ws.Range("H61").Value = 2486.8572
$ws.Range("I61").Value = 2541.4
$ws.Range("J61").Value = 2456.5557
$ws.Range("K61").Value = 2541.4
$ws.Range("L61").Value = 2456.5557
$ws.Range("M61").Value = -2329.4
$ws.Range("N61").Value = -2880.5557
$ws.Range("H74").Value = 56392.25
$ws.Range("I74").Value = 7325
$ws.Range("K74").Value = 7325
$ws.Range("M74").Value = -6451
$ws.Range("H77").Value = 56392.25
$ws.Range("I77").Value = 7325
$ws.Range("K77").Value = 36625
$ws.Range("M77").Value = -32257
$ws.Range("H122").Value = 11699522
$ws.Range("I122").Value = 18521460
$ws.Range("K122").Value = 55564380
$ws.Range("M122").Value = -55561930
$ws.Range("H132").Value = 2271
$ws.Range("I132").Value = 1269.0714
$ws.Range("J132").Value = 3272.9285
$ws.Range("K132").Value = 3807.2142
$ws.Range("L132").Value = 9818.7855
$ws.Range("M132").Value = -1277.2142
$ws.Range("N132").Value = -14878.7855
$ws.Range("H136").Value = 2486.8572
$ws.Range("I136").Value = 2541.4
$ws.Range("J136").Value = 2456.5557
$ws.Range("K136").Value = 7624.200000000001
$ws.Range("L136").Value = 7369.6671
$ws.Range("M136").Value = -5074.200000000001
$ws.Range("N136").Value = -12469.6671
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1922.6
$ws.Range("I16").Value = 1647.6666
$ws.Range("K16").Value = 1647.6666
$ws.Range("M16").Value = -1360.6666
$ws.Range("H31").Value = 14688.153
$ws.Range("J31").Value = 18291.016
$ws.Range("L31").Value = 18291.016
$ws.Range("N31").Value = -18881.016
$ws.Range("H34").Value = 14688.153
$ws.Range("J34").Value = 18291.016
$ws.Range("L34").Value = 18291.016
$ws.Range("N34").Value = -18695.016
$ws.Range("H62").Value = 2242.7144
$ws.Range("I62").Value = 1233.3334
$ws.Range("J62").Value = 2999.75
$ws.Range("K62").Value = 1233.3334
$ws.Range("L62").Value = 2999.75
$ws.Range("M62").Value = -609.3334
$ws.Range("N62").Value = -4247.75
$ws.Range("H65").Value = 2242.7144
$ws.Range("I65").Value = 1233.3334
$ws.Range("J65").Value = 2999.75
$ws.Range("K65").Value = 6166.666999999999
$ws.Range("L65").Value = 14998.75
$ws.Range("M65").Value = -3046.666999999999
$ws.Range("N65").Value = -21238.75
$ws.Range("H99").Value = 3502.5715
$ws.Range("I99").Value = 2724.5715
$ws.Range("J99").Value = 4280.5713
$ws.Range("K99").Value = 2724.5715
$ws.Range("L99").Value = 4280.5713
$ws.Range("M99").Value = -1226.5715
$ws.Range("N99").Value = -7276.5713
$ws.Range("H113").Value = 1922.6
$ws.Range("I113").Value = 1647.6666
$ws.Range("K113").Value = 1647.6666
$ws.Range("M113").Value = 522.3334
$ws.Range("H126").Value = 3502.5715
$ws.Range("I126").Value = 2724.5715
$ws.Range("J126").Value = 4280.5713
$ws.Range("K126").Value = 8173.7145
$ws.Range("L126").Value = 12841.7139
$ws.Range("M126").Value = -5703.7145
$ws.Range("N126").Value = -17781.7139
$ws.Range("H132").Value = 116484.5
$ws.Range("I132").Value = 5998.1665
$ws.Range("J132").Value = 447943.5
$ws.Range("K132").Value = 17994.4995
$ws.Range("L132").Value = 1343830.5
$ws.Range("M132").Value = -15464.4995
$ws.Range("N132").Value = -1348890.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 7078.9
$ws.Range("I94").Value = 3531.3333
$ws.Range("K94").Value = 10593.9999
$ws.Range("M94").Value = -9917.999899999999
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("N112").ClearContents()
$ws.Range("H131").Value = 19844656
$ws.Range("I131").Value = 20834032
$ws.Range("J131").Value = 19611862
$ws.Range("K131").Value = 62502096
$ws.Range("L131").Value = 58835586
$ws.Range("M131").Value = -62497056
$ws.Range("N131").Value = -58845666
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8516.083000000001
$ws.Range("I70").Value = 9389.556
$ws.Range("J70").Value = 5895.6665
$ws.Range("K70").Value = 9389.556
$ws.Range("L70").Value = 5895.6665
$ws.Range("M70").Value = -9119.556
$ws.Range("N70").Value = -6435.6665
$ws.Range("H73").Value = 8516.083000000001
$ws.Range("I73").Value = 9389.556
$ws.Range("J73").Value = 5895.6665
$ws.Range("K73").Value = 9389.556
$ws.Range("L73").Value = 5895.6665
$ws.Range("M73").Value = -8453.556
$ws.Range("N73").Value = -7767.6665
$ws.Range("H102").Value = 29856.244
$ws.Range("I102").Value = 2367.88
$ws.Range("J102").Value = 87123.664
$ws.Range("K102").Value = 2367.88
$ws.Range("L102").Value = 87123.664
$ws.Range("M102").Value = -745.8800000000001
$ws.Range("N102").Value = -90367.664
$ws.Range("H107").Value = 779.2
$ws.Range("I107").Value = 779.2
$ws.Range("K107").Value = 779.2
$ws.Range("M107").Value = 1140.8
$ws.Range("H122").Value = 205669.31
$ws.Range("I122").Value = 298901.38
$ws.Range("J122").Value = 5886.357
$ws.Range("K122").Value = 896704.14
$ws.Range("L122").Value = 17659.071
$ws.Range("M122").Value = -894254.14
$ws.Range("N122").Value = -22559.071
$ws.Range("H126").Value = 3619.2666
$ws.Range("I126").Value = 3382.9167
$ws.Range("K126").Value = 10148.7501
$ws.Range("M126").Value = -7678.750100000001
$ws.Range("H132").Value = 4901.9165
$ws.Range("I132").Value = 3992.5
$ws.Range("K132").Value = 11977.5
$ws.Range("M132").Value = -9447.5
$ws.Range("H136").Value = 58498
$ws.Range("I136").Value = 30000
$ws.Range("J136").Value = 65622.5
$ws.Range("K136").Value = 90000
$ws.Range("L136").Value = 196867.5
$ws.Range("M136").Value = -87450
$ws.Range("N136").Value = -201967.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H136").Value = 2573.8965
$ws.Range("I136").Value = 2131.0417
$ws.Range("K136").Value = 6393.125100000001
$ws.Range("M136").Value = -3843.125100000001
